# Scheduled runner update: refresh Universalis market-price pulls
# (currentAveragePrice* / LevePrice* / LeveProfit* columns) across the
# per-job Leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ALC!row11
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 626.5714
$ws.Range("I11").Value = 626.5714
$ws.Range("K11").Value = 626.5714
$ws.Range("M11").Value = -486.5714

# ALC!row43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 6661.25
$ws.Range("I43").Value = 5495.5
$ws.Range("K43").Value = 5495.5
$ws.Range("M43").Value = -5426.5

# ALC!row62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3315.75
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 3421
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 3421
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -4669

# ALC!row65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 3315.75
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 3421
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 17105
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -23345

# ALC!row107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 57020.5
$ws.Range("J107").Value = 1250
$ws.Range("L107").Value = 1250
$ws.Range("N107").Value = -5090

# ALC!row113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

# ALC!row137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3019.6
$ws.Range("J137").Value = 4366
$ws.Range("L137").Value = 13098
$ws.Range("N137").Value = -18198

# ARM!row2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 366.33334
$ws.Range("I2").Value = 437
$ws.Range("J2").Value = 225
$ws.Range("K2").Value = 437
$ws.Range("L2").Value = 225
$ws.Range("M2").Value = -324
$ws.Range("N2").Value = -451

# ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1004.0909
$ws.Range("I32").Value = 1004.0909
$ws.Range("K32").Value = 1004.0909
$ws.Range("M32").Value = -717.0909

# ARM!row61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2720.4
$ws.Range("I61").Value = 2529.3333
$ws.Range("J61").Value = 3007
$ws.Range("K61").Value = 2529.3333
$ws.Range("L61").Value = 3007
$ws.Range("M61").Value = -2317.3333
$ws.Range("N61").Value = -3431

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1555.75
$ws.Range("I74").Value = 1111.5
$ws.Range("K74").Value = 1111.5
$ws.Range("M74").Value = -237.5

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1555.75
$ws.Range("I77").Value = 1111.5
$ws.Range("K77").Value = 5557.5
$ws.Range("M77").Value = -1189.5

# ARM!row102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 42200588
$ws.Range("I102").Value = 3666814.8
$ws.Range("J102").Value = 100001250
$ws.Range("K102").Value = 3666814.8
$ws.Range("L102").Value = 100001250
$ws.Range("M102").Value = -3665192.8
$ws.Range("N102").Value = -100004494

# ARM!row116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 366.33334
$ws.Range("I116").Value = 437
$ws.Range("J116").Value = 225
$ws.Range("K116").Value = 437
$ws.Range("L116").Value = 225
$ws.Range("M116").Value = 1857
$ws.Range("N116").Value = -4813

# ARM!row136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2720.4
$ws.Range("I136").Value = 2529.3333
$ws.Range("J136").Value = 3007
$ws.Range("K136").Value = 7587.999899999999
$ws.Range("L136").Value = 9021
$ws.Range("M136").Value = -5037.999899999999
$ws.Range("N136").Value = -14121

# BSM!row3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 366.33334
$ws.Range("I3").Value = 437
$ws.Range("J3").Value = 225
$ws.Range("K3").Value = 437
$ws.Range("L3").Value = 225
$ws.Range("M3").Value = -323
$ws.Range("N3").Value = -453

# BSM!row94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 65866.7
$ws.Range("I94").Value = 79484.57000000001
$ws.Range("J94").Value = 2316.6667
$ws.Range("K94").Value = 79484.57000000001
$ws.Range("L94").Value = 2316.6667
$ws.Range("M94").Value = -79033.57000000001
$ws.Range("N94").Value = -3218.6667

# BSM!row99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 5010
$ws.Range("I99").Value = 5010
$ws.Range("K99").Value = 5010
$ws.Range("M99").Value = -3512

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5186.933
$ws.Range("J31").Value = 7268.857
$ws.Range("L31").Value = 7268.857
$ws.Range("N31").Value = -7858.857

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5186.933
$ws.Range("J34").Value = 7268.857
$ws.Range("L34").Value = 7268.857
$ws.Range("N34").Value = -7672.857

# CRP!row58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2006.5
$ws.Range("J58").Value = 3000
$ws.Range("L58").Value = 3000
$ws.Range("N58").Value = -3406

# CRP!row99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 501587.8
$ws.Range("I99").Value = 1800.25
$ws.Range("J99").Value = 2500738
$ws.Range("K99").Value = 1800.25
$ws.Range("L99").Value = 2500738
$ws.Range("M99").Value = -302.25
$ws.Range("N99").Value = -2503734

# CRP!row126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 501587.8
$ws.Range("I126").Value = 1800.25
$ws.Range("J126").Value = 2500738
$ws.Range("K126").Value = 5400.75
$ws.Range("L126").Value = 7502214
$ws.Range("M126").Value = -2930.75
$ws.Range("N126").Value = -7507154

# CRP!row136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2006.5
$ws.Range("J136").Value = 3000
$ws.Range("L136").Value = 9000
$ws.Range("N136").Value = -14100

# CRP!row140
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# CUL!row7
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 233.33333
$ws.Range("I7").Value = 100
$ws.Range("K7").Value = 300
$ws.Range("M7").Value = -188

# CUL!row132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 750
$ws.Range("J132").Value = 1000
$ws.Range("L132").Value = 9000
$ws.Range("N132").Value = -14060

# GSM!row102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1643.5
$ws.Range("I102").Value = 1372.2
$ws.Range("K102").Value = 1372.2
$ws.Range("M102").Value = 249.8

# LTW!row132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2432
$ws.Range("I132").Value = 1252
$ws.Range("J132").Value = 2904
$ws.Range("K132").Value = 3756
$ws.Range("L132").Value = 8712
$ws.Range("M132").Value = -1226
$ws.Range("N132").Value = -13772

# LTW!row136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3502
$ws.Range("I136").Value = 3502
$ws.Range("K136").Value = 10506
$ws.Range("M136").Value = -7956

# WVR!row126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3870.182
$ws.Range("I126").Value = 3256.7
$ws.Range("K126").Value = 9770.099999999999
$ws.Range("M126").Value = -7300.099999999999

